$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 175. This shifts the previous rows
# 175:279 down to 176:280 (so the sheet grows from 279 to 280 data rows,
# dimension A1:T279 -> A1:T280).
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with its data. The static /
# categorical columns (A,B,C,E,F,G,H,I,J,K,Q,T) are identical for every
# row in this sheet, and L/R keep the same values that were already in
# the (now shifted) row, while D/M/N/O/P/S get the new reported values.
$ws.Cells.Item(175, 1).Value  = 10
$ws.Cells.Item(175, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(175, 3).Value  = "La Araucanía"
$ws.Cells.Item(175, 4).Value  = 44572
$ws.Cells.Item(175, 5).Value  = 9
$ws.Cells.Item(175, 6).Value  = "Fruta"
$ws.Cells.Item(175, 7).Value  = 100108
$ws.Cells.Item(175, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(175, 9).Value  = 100108002
$ws.Cells.Item(175, 10).Value = "Mango"
$ws.Cells.Item(175, 11).Value = "Sin especificar"
$ws.Cells.Item(175, 12).Value = "Primera"
$ws.Cells.Item(175, 13).Value = 800
$ws.Cells.Item(175, 14).Value = 7000
$ws.Cells.Item(175, 15).Value = 7000
$ws.Cells.Item(175, 16).Value = 7000
$ws.Cells.Item(175, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(175, 18).Value = "Perú"
$ws.Cells.Item(175, 19).Value = 1750
$ws.Cells.Item(175, 20).Value = 4
